$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1436925.2
$ws.Range("I11").Value = 1436925.2
$ws.Range("K11").Value = 1436925.2
$ws.Range("M11").Value = -1436785.2
$ws.Range("H32").Value = 257.15384
$ws.Range("I32").Value = 272.16666
$ws.Range("J32").Value = 244.28572
$ws.Range("K32").Value = 272.16666
$ws.Range("L32").Value = 244.28572
$ws.Range("M32").Value = 53.83334000000002
$ws.Range("N32").Value = -896.28572
$ws.Range("H55").Value = 163
$ws.Range("J55").Value = 198
$ws.Range("L55").Value = 198
$ws.Range("N55").Value = -626
$ws.Range("H107").Value = 1914
$ws.Range("J107").Value = 750
$ws.Range("L107").Value = 750
$ws.Range("N107").Value = -4590
$ws.Range("H132").Value = 102325.15
$ws.Range("I132").Value = 119159.15
$ws.Range("K132").Value = 357477.45
$ws.Range("M132").Value = -354947.45
$ws.Range("H141").Value = 3569
$ws.Range("I141").Value = 3373.2222
$ws.Range("J141").Value = 4450
$ws.Range("K141").Value = 10119.6666
$ws.Range("L141").Value = 13350
$ws.Range("M141").Value = -4939.6666
$ws.Range("N141").Value = -23710
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6470.2856
$ws.Range("I32").Value = 4568.7026
$ws.Range("J32").Value = 10173.368
$ws.Range("K32").Value = 4568.7026
$ws.Range("L32").Value = 10173.368
$ws.Range("M32").Value = -4281.7026
$ws.Range("N32").Value = -10747.368
$ws.Range("H45").Value = 1187.5
$ws.Range("I45").Value = 1050
$ws.Range("J45").Value = 1233.3334
$ws.Range("K45").Value = 1050
$ws.Range("L45").Value = 1233.3334
$ws.Range("M45").Value = -673
$ws.Range("N45").Value = -1987.3334
$ws.Range("H61").Value = 2074.44
$ws.Range("I61").Value = 1798.2273
$ws.Range("K61").Value = 1798.2273
$ws.Range("M61").Value = -1586.2273
$ws.Range("H131").Value = 41715
$ws.Range("J131").Value = 41715
$ws.Range("L131").Value = 41715
$ws.Range("N131").Value = -51795
$ws.Range("H132").Value = 2598.653
$ws.Range("I132").Value = 2000.9
$ws.Range("K132").Value = 6002.700000000001
$ws.Range("M132").Value = -3472.700000000001
$ws.Range("H136").Value = 2074.44
$ws.Range("I136").Value = 1798.2273
$ws.Range("K136").Value = 5394.6819
$ws.Range("M136").Value = -2844.6819
$ws.Range("H137").Value = 40621
$ws.Range("J137").Value = 40621
$ws.Range("L137").Value = 40621
$ws.Range("N137").Value = -50821
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 208.51852
$ws.Range("I80").Value = 98.666664
$ws.Range("J80").Value = 263.44446
$ws.Range("K80").Value = 98.666664
$ws.Range("L80").Value = 263.44446
$ws.Range("M80").Value = 899.333336
$ws.Range("N80").Value = -2259.44446
$ws.Range("H83").Value = 208.51852
$ws.Range("I83").Value = 98.666664
$ws.Range("J83").Value = 263.44446
$ws.Range("K83").Value = 493.33332
$ws.Range("L83").Value = 1317.2223
$ws.Range("M83").Value = 4498.66668
$ws.Range("N83").Value = -11301.2223
$ws.Range("H105").Value = 2575.1428
$ws.Range("I105").Value = 2534.0625
$ws.Range("J105").Value = 3013.3333
$ws.Range("K105").Value = 2534.0625
$ws.Range("L105").Value = 3013.3333
$ws.Range("M105").Value = -787.0625
$ws.Range("N105").Value = -6507.3333
$ws.Range("H107").Value = 1043.6428
$ws.Range("I107").Value = 947.619
$ws.Range("J107").Value = 1331.7142
$ws.Range("K107").Value = 947.619
$ws.Range("L107").Value = 1331.7142
$ws.Range("M107").Value = 972.381
$ws.Range("N107").Value = -5171.7142
$ws.Range("H134").Value = 2566.8406
$ws.Range("I134").Value = 1594.5807
$ws.Range("K134").Value = 4783.742099999999
$ws.Range("M134").Value = -2248.742099999999
$ws.Range("H137").Value = 39490
$ws.Range("J137").Value = 40544.445
$ws.Range("L137").Value = 40544.445
$ws.Range("N137").Value = -50744.445
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11630919
$ws.Range("I31").Value = 1614.7241
$ws.Range("K31").Value = 1614.7241
$ws.Range("M31").Value = -1319.7241
$ws.Range("H34").Value = 11630919
$ws.Range("I34").Value = 1614.7241
$ws.Range("K34").Value = 1614.7241
$ws.Range("M34").Value = -1412.7241
$ws.Range("H58").Value = 1860.403
$ws.Range("I58").Value = 1629.1608
$ws.Range("J58").Value = 3037.6365
$ws.Range("K58").Value = 1629.1608
$ws.Range("L58").Value = 3037.6365
$ws.Range("M58").Value = -1426.1608
$ws.Range("N58").Value = -3443.6365
$ws.Range("H99").Value = 9527001
$ws.Range("I99").Value = 20002028
$ws.Range("J99").Value = 4249
$ws.Range("K99").Value = 20002028
$ws.Range("L99").Value = 4249
$ws.Range("M99").Value = -20000530
$ws.Range("N99").Value = -7245
$ws.Range("H126").Value = 9527001
$ws.Range("I126").Value = 20002028
$ws.Range("J126").Value = 4249
$ws.Range("K126").Value = 60006084
$ws.Range("L126").Value = 12747
$ws.Range("M126").Value = -60003614
$ws.Range("N126").Value = -17687
$ws.Range("H132").Value = 2912.9697
$ws.Range("I132").Value = 1281.9375
$ws.Range("J132").Value = 4448.0586
$ws.Range("K132").Value = 3845.8125
$ws.Range("L132").Value = 13344.1758
$ws.Range("M132").Value = -1315.8125
$ws.Range("N132").Value = -18404.1758
$ws.Range("H134").Value = 4833.4287
$ws.Range("I134").Value = 5022.64
$ws.Range("J134").Value = 4360.4
$ws.Range("K134").Value = 15067.92
$ws.Range("L134").Value = 13081.2
$ws.Range("M134").Value = -12532.92
$ws.Range("N134").Value = -18151.2
$ws.Range("H136").Value = 1860.403
$ws.Range("I136").Value = 1629.1608
$ws.Range("J136").Value = 3037.6365
$ws.Range("K136").Value = 4887.482400000001
$ws.Range("L136").Value = 9112.9095
$ws.Range("M136").Value = -2337.482400000001
$ws.Range("N136").Value = -14212.9095
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 928698.1
$ws.Range("I4").Value = 8036518
$ws.Range("J4").Value = 1591.174
$ws.Range("K4").Value = 24109554
$ws.Range("L4").Value = 4773.522
$ws.Range("M4").Value = -24109442
$ws.Range("N4").Value = -4997.522
$ws.Range("H5").Value = 1958.7894
$ws.Range("I5").Value = 396.1111
$ws.Range("J5").Value = 3365.2
$ws.Range("K5").Value = 1188.3333
$ws.Range("L5").Value = 10095.6
$ws.Range("M5").Value = -1076.3333
$ws.Range("N5").Value = -10319.6
$ws.Range("H86").Value = 5251.5
$ws.Range("J86").Value = 6835.3335
$ws.Range("L86").Value = 20506.0005
$ws.Range("N86").Value = -22878.0005
$ws.Range("H89").Value = 5251.5
$ws.Range("J89").Value = 6835.3335
$ws.Range("L89").Value = 61518.0015
$ws.Range("N89").Value = -73374.0015
$ws.Range("H131").Value = 8621655
$ws.Range("J131").Value = 861.8431399999999
$ws.Range("L131").Value = 2585.52942
$ws.Range("N131").Value = -12665.52942
$ws.Range("H135").Value = 1958.7894
$ws.Range("I135").Value = 396.1111
$ws.Range("J135").Value = 3365.2
$ws.Range("K135").Value = 3564.9999
$ws.Range("L135").Value = 30286.8
$ws.Range("M135").Value = -1029.9999
$ws.Range("N135").Value = -35356.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 52852.777
$ws.Range("J137").Value = 58890
$ws.Range("L137").Value = 58890
$ws.Range("N137").Value = -69090
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6899.816
$ws.Range("I40").Value = 5447.857
$ws.Range("J40").Value = 8693.412
$ws.Range("K40").Value = 5447.857
$ws.Range("L40").Value = 8693.412
$ws.Range("M40").Value = -5311.857
$ws.Range("N40").Value = -8965.412
$ws.Range("H55").Value = 317.17648
$ws.Range("I55").Value = 200.72728
$ws.Range("J55").Value = 530.6667
$ws.Range("K55").Value = 200.72728
$ws.Range("L55").Value = 530.6667
$ws.Range("M55").Value = -27.72728000000001
$ws.Range("N55").Value = -876.6667
$ws.Range("H122").Value = 6706.1177
$ws.Range("I122").Value = 4286.2856
$ws.Range("J122").Value = 8400
$ws.Range("K122").Value = 12858.8568
$ws.Range("L122").Value = 25200
$ws.Range("M122").Value = -10408.8568
$ws.Range("N122").Value = -30100
$ws.Range("H136").Value = 4148.1904
$ws.Range("I136").Value = 2186.5715
$ws.Range("J136").Value = 8071.4287
$ws.Range("K136").Value = 6559.7145
$ws.Range("L136").Value = 24214.2861
$ws.Range("M136").Value = -4009.7145
$ws.Range("N136").Value = -29314.2861
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11496339
$ws.Range("I132").Value = 1105.75
$ws.Range("J132").Value = 25644320
$ws.Range("K132").Value = 3317.25
$ws.Range("L132").Value = 76932960
$ws.Range("M132").Value = -787.25
$ws.Range("N132").Value = -76938020
$ws.Range("H136").Value = 1857.8
$ws.Range("I136").Value = 705.6
$ws.Range("K136").Value = 2116.8
$ws.Range("M136").Value = 433.1999999999998
